# Update workbook to xsd 1.5.9
# - Readme sheet: bump version/date strings
# - Admin sheet: update lookup lists (Library Strategy, Instrument, File Type)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Readme sheet: update the "last updated" banner and version number
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("Readme")
$readme.Range("A1").Value = "2021-07-13 Bioinformation and DDBJ Center"
$readme.Range("A2").Value = "v1.2"

# ---------------------------------------------------------------------------
# Admin sheet: refresh the lookup/reference tables used for data validation
# ---------------------------------------------------------------------------
$admin = $wb.Worksheets.Item("Admin")

# Column C - Library Strategy (NOMe-Seq inserted before RIP-Seq)
$CVals = @("Library Strategy","WGS","WGA","WXS","RNA-Seq","ssRNA-seq","miRNA-Seq","ncRNA-Seq","FL-cDNA","EST","Hi-C","ATAC-seq","WCS","RAD-Seq","CLONE","POOLCLONE","AMPLICON","CLONEEND","FINISHING","ChIP-Seq","MNase-Seq","DNase-Hypersensitivity","Bisulfite-Seq","CTS","MRE-Seq","MeDIP-Seq","MBD-Seq","Tn-Seq","VALIDATION","FAIRE-seq","SELEX","NOMe-Seq","RIP-Seq","ChIA-PET","Synthetic-Long-Read","Targeted-Capture","Tethered Chromatin Conformation Capture","OTHER")

# Column D - Instrument (trimmed names + many new sequencer models)
$DVals = @("Instrument","454 GS","454 GS 20","454 GS FLX","454 GS FLX+","454 GS FLX Titanium","454 GS Junior","Illumina Genome Analyzer","Illumina Genome Analyzer II","Illumina Genome Analyzer IIx","Illumina HiSeq 1000","Illumina HiSeq 1500","Illumina HiSeq 2000","Illumina HiSeq 2500","Illumina HiSeq 3000","Illumina HiSeq 4000","HiSeq X Five","HiSeq X Ten","Illumina HiScanSQ","Illumina NovaSeq 6000","Illumina MiSeq","Illumina MiniSeq","Illumina iSeq 100","NextSeq 500","NextSeq 550","NextSeq 1000","NextSeq 2000","Helicos HeliScope","AB SOLiD System","AB SOLiD System 2.0","AB SOLiD System 3.0","AB SOLiD 3 Plus System","AB SOLiD 4 System","AB SOLiD 4hq System","AB SOLiD PI System","AB 5500 Genetic Analyzer","AB 5500xl Genetic Analyzer","AB 5500xl-W Genetic Analysis System","Complete Genomics","BGISEQ-500","DNBSEQ-G400","DNBSEQ-T7","DNBSEQ-G50","MGISEQ-2000RS","PacBio RS","PacBio RS II","Sequel","Sequel II","Ion Torrent PGM","Ion Torrent Proton","Ion Torrent S5","Ion Torrent S5 XL","Ion GeneStudio S5","Ion GeneStudio S5 plus","Ion GeneStudio S5 prime","AB 3730xL Genetic Analyzer","AB 3730 Genetic Analyzer","AB 3500xL Genetic Analyzer","AB 3500 Genetic Analyzer","AB 3130xL Genetic Analyzer","AB 3130 Genetic Analyzer","AB 310 Genetic Analyzer","MinION","GridION","PromethION","unspecified")

# Column F - File Type ("fastq" entry removed, list shifts up)
$FVals = @("File Type","generic_fastq","sff","PacBio_HDF5","bam","tab","reference_fasta","SOLiD_native","")

for ($i = 0; $i -lt $CVals.Length; $i++) {
    $admin.Cells.Item($i + 1, 3).Value = $CVals[$i]
}

for ($i = 0; $i -lt $DVals.Length; $i++) {
    $admin.Cells.Item($i + 1, 4).Value = $DVals[$i]
}

for ($i = 0; $i -lt $FVals.Length; $i++) {
    $admin.Cells.Item($i + 1, 6).Value = $FVals[$i]
}
